{"js": "// Replace the three-digit \u00f7 one-digit division answers in the table\n// with the updated values, one run of text per table cell.\nconst replacements = [\n  [\"728\u00f72=364, 0\", \"512\u00f74=128, 0\"],\n  [\"268\u00f77=38, 2\", \"898\u00f77=128, 2\"],\n  [\"568\u00f79=63, 1\", \"393\u00f79=43, 6\"],\n  [\"503\u00f76=83, 5\", \"789\u00f72=394, 1\"],\n  [\"645\u00f74=161, 1\", \"706\u00f78=88, 2\"],\n  [\"237\u00f76=39, 3\", \"801\u00f75=160, 1\"],\n  [\"580\u00f78=72, 4\", \"527\u00f76=87, 5\"],\n  [\"870\u00f72=435, 0\", \"437\u00f72=218, 1\"],\n  [\"820\u00f72=410, 0\", \"977\u00f75=195, 2\"],\n  [\"943\u00f75=188, 3\", \"122\u00f74=30, 2\"],\n  [\"955\u00f74=238, 3\", \"595\u00f76=99, 1\"],\n  [\"591\u00f79=65, 6\", \"704\u00f73=234, 2\"],\n  [\"457\u00f76=76, 1\", \"489\u00f72=244, 1\"],\n  [\"608\u00f74=152, 0\", \"933\u00f76=155, 3\"],\n  [\"148\u00f76=24, 4\", \"390\u00f75=78, 0\"],\n  [\"408\u00f77=58, 2\", \"690\u00f76=115, 0\"],\n  [\"756\u00f77=108, 0\", \"105\u00f72=52, 1\"],\n  [\"314\u00f75=62, 4\", \"653\u00f73=217, 2\"],\n  [\"456\u00f79=50, 6\", \"267\u00f73=89, 0\"],\n  [\"841\u00f74=210, 1\", \"461\u00f79=51, 2\"],\n  [\"326\u00f73=108, 2\", \"658\u00f77=94, 0\"],\n  [\"138\u00f74=34, 2\", \"922\u00f74=230, 2\"],\n  [\"995\u00f75=199, 0\", \"756\u00f79=84, 0\"],\n  [\"489\u00f76=81, 3\", \"587\u00f74=146, 3\"],\n  [\"891\u00f78=111, 3\", \"614\u00f75=122, 4\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the three-digit / one-digit division answers in the table\n# with the updated values, using Word's Find/Replace on each cell's text.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @('728\u00f72=364, 0', '512\u00f74=128, 0'),\n    @('268\u00f77=38, 2', '898\u00f77=128, 2'),\n    @('568\u00f79=63, 1', '393\u00f79=43, 6'),\n    @('503\u00f76=83, 5', '789\u00f72=394, 1'),\n    @('645\u00f74=161, 1', '706\u00f78=88, 2'),\n    @('237\u00f76=39, 3', '801\u00f75=160, 1'),\n    @('580\u00f78=72, 4', '527\u00f76=87, 5'),\n    @('870\u00f72=435, 0', '437\u00f72=218, 1'),\n    @('820\u00f72=410, 0', '977\u00f75=195, 2'),\n    @('943\u00f75=188, 3', '122\u00f74=30, 2'),\n    @('955\u00f74=238, 3', '595\u00f76=99, 1'),\n    @('591\u00f79=65, 6', '704\u00f73=234, 2'),\n    @('457\u00f76=76, 1', '489\u00f72=244, 1'),\n    @('608\u00f74=152, 0', '933\u00f76=155, 3'),\n    @('148\u00f76=24, 4', '390\u00f75=78, 0'),\n    @('408\u00f77=58, 2', '690\u00f76=115, 0'),\n    @('756\u00f77=108, 0', '105\u00f72=52, 1'),\n    @('314\u00f75=62, 4', '653\u00f73=217, 2'),\n    @('456\u00f79=50, 6', '267\u00f73=89, 0'),\n    @('841\u00f74=210, 1', '461\u00f79=51, 2'),\n    @('326\u00f73=108, 2', '658\u00f77=94, 0'),\n    @('138\u00f74=34, 2', '922\u00f74=230, 2'),\n    @('995\u00f75=199, 0', '756\u00f79=84, 0'),\n    @('489\u00f76=81, 3', '587\u00f74=146, 3'),\n    @('891\u00f78=111, 3', '614\u00f75=122, 4'),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $found = $find.Execute($null, $true, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
